$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds text-formatted numbers (e.g. "51.861.28", "3.40").
# Force Text format cell-by-cell before writing so Excel does not silently
# reinterpret numeric-looking strings as numbers (which would drop trailing
# zeros / change thousands-separator formatting).
$priceCells = @(
    'D2',
    'D3',
    'D5',
    'D6',
    'D10',
    'D13',
    'D14',
    'D15',
    'D17',
    'D18',
    'D19',
    'D20',
    'D21',
    'D22',
    'D23',
    'D24',
    'D25',
    'D26',
    'D28',
    'D29',
    'D30',
    'D32',
    'D33',
    'D34',
    'D35',
    'D36',
    'D38',
    'D42',
    'D43',
    'D44',
    'D45',
    'D46',
    'D47',
    'D48',
    'D49',
    'D50',
    'D51'
)
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range('D2').Value = '51.861.28'
$ws.Range('D3').Value = '2.844.83'
$ws.Range('E3').Value = '  +1.92%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '350.48'
$ws.Range('E5').Value = '  -1.08%  '
$ws.Range('D6').Value = '113.02'
$ws.Range('E6').Value = '  +3.36%  '
$ws.Range('E7').Value = '  +1.22%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +3.60%  '
$ws.Range('D10').Value = '40.22'
$ws.Range('E10').Value = '  +0.57%  '
$ws.Range('E11').Value = '  -0.91%  '
$ws.Range('E12').Value = '  +1.05%  '
$ws.Range('D13').Value = '20.13'
$ws.Range('D14').Value = '7.79'
$ws.Range('E14').Value = '  +1.90%  '
$ws.Range('D15').Value = '3.278.23'
$ws.Range('E15').Value = '  +1.38%  '
$ws.Range('E16').Value = '  +6.26%  '
$ws.Range('D17').Value = '2.852.05'
$ws.Range('E17').Value = '  +1.20%  '
$ws.Range('D18').Value = '51.924.14'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('D19').Value = '3.40'
$ws.Range('E19').Value = '  +9.09%  '
$ws.Range('D20').Value = '7.65'
$ws.Range('E20').Value = '  -1.21%  '
$ws.Range('D21').Value = '13.45'
$ws.Range('E21').Value = '  +2.19%  '
$ws.Range('D22').Value = '0.0₃0973'
$ws.Range('E22').Value = '  +0.66%  '
$ws.Range('D23').Value = '70.46'
$ws.Range('E23').Value = '  +0.57%  '
$ws.Range('D24').Value = '268.88'
$ws.Range('E24').Value = '  +0.87%  '
$ws.Range('D25').Value = '2.75'
$ws.Range('E25').Value = '  +0.78%  '
$ws.Range('D26').Value = '26.31'
$ws.Range('E26').Value = '  +0.48%  '
$ws.Range('E27').Value = '  +0.10%  '
$ws.Range('D28').Value = '0.163'
$ws.Range('E28').Value = '  +0.70%  '
$ws.Range('D29').Value = '39.42'
$ws.Range('E29').Value = '  +6.62%  '
$ws.Range('D30').Value = '10.54'
$ws.Range('E30').Value = '  +2.90%  '
$ws.Range('E31').Value = '  +15.92%  '
$ws.Range('D32').Value = '6.29'
$ws.Range('E32').Value = '  +0.82%  '
$ws.Range('D33').Value = '52.76'
$ws.Range('E33').Value = '  +1.49%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.0894'
$ws.Range('E34').Value = '  +7.61%  '
$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').Value = '0.0450'
$ws.Range('E35').Value = '  -1.18%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').Value = '5.62'
$ws.Range('E36').Value = '  +0.44%  '
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('D38').Value = '18.93'
$ws.Range('E38').Value = '  +1.95%  '
$ws.Range('E39').Value = '  +2.45%  '
$ws.Range('E40').Value = '  +2.24%  '
$ws.Range('E41').Value = '  +1.42%  '
$ws.Range('D42').Value = '2.53'
$ws.Range('E42').Value = '  -1.02%  '
$ws.Range('D43').Value = '122.99'
$ws.Range('E43').Value = '  +1.68%  '
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = '2.23'
$ws.Range('E44').Value = '  +1.91%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '22.21'
$ws.Range('E45').Value = '  +0.75%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = '2.52'
$ws.Range('E46').Value = '  +8.22%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = '3.52'
$ws.Range('E47').Value = '  +6.86%  '
$ws.Range('D48').Value = '2.175.87'
$ws.Range('E48').Value = '  +1.80%  '
$ws.Range('D49').Value = '0.248'
$ws.Range('E49').Value = '  +21.22%  '
$ws.Range('D50').Value = '0.950'
$ws.Range('E50').Value = '  +4.38%  '
$ws.Range('D51').Value = '5.43'
$ws.Range('E51').Value = '  -0.44%  '
